# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" (countries) ranking table with the latest
# figures and re-sort the affected rows by total cases (column B, desc).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-Row([int]$row, [string]$pais, $casos, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value2 = $pais
    $ws.Cells.Item($row, 2).Value2 = $casos
    $ws.Cells.Item($row, 3).Value2 = $nuevos
    $ws.Cells.Item($row, 4).Value2 = $activos
    $ws.Cells.Item($row, 5).Value2 = $recuperados
    $ws.Cells.Item($row, 6).Value2 = $criticos
    $ws.Cells.Item($row, 7).Value2 = $muertesHoy
    $ws.Cells.Item($row, 8).Value2 = $muertes
}

# Timestamp banner
$ws.Range("A1").Value2 = "Datos actualizados a 4 de Abril de 2020 a las 21:52"

# --- Estados Unidos (row 4): updated totals, stays #1 ---
Set-Row 4 "Estados Unidos" 304383 27222 14686 281421 7983 872 8276

# --- Suiza (row 13): active/recovered split updated, total unchanged ---
Set-Row 13 "Suiza" 20505 899 6415 13424 391 75 666

# --- Brasil overtakes Corea del Sur (rows 19-20) ---
Set-Row 19 "Brasil" 10278 1084 127 9720 296 68 431
Set-Row 20 "Corea del Sur" 10156 94 6325 3654 55 3 177

# --- Sudafrica (row 48): updated totals ---
Set-Row 48 "Sudafrica" 1585 80 95 1481 7 0 9

# --- Colombia overtakes Argentina and Catar (rows 52-54) ---
Set-Row 52 "Colombia" 1406 139 85 1289 50 7 32
Set-Row 53 "Argentina" 1353 0 279 1032 0 0 42
Set-Row 54 "Catar" 1325 250 109 1213 37 0 3

# --- Tunez (row 73): minor update ---
Set-Row 73 "Tunez" 553 58 5 530 26 0 18

# --- Monaco overtakes Aruba (rows 131-132) ---
Set-Row 131 "Monaco" 66 2 3 62 2 0 1
Set-Row 132 "Aruba" 64 2 1 63 0 0 0

# --- Curazao (row 176): minor update ---
Set-Row 176 "Curazao" 11 0 5 5 0 0 1
